$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Timestamp (column A) values for rows 2..97
$aValues = @(
46079.01041666666,
46079.02083333334,
46079.03125,
46079.04166666666,
46079.05208333334,
46079.0625,
46079.07291666666,
46079.08333333334,
46079.09375,
46079.10416666666,
46079.11458333334,
46079.125,
46079.13541666666,
46079.14583333334,
46079.15625,
46079.16666666666,
46079.17708333334,
46079.1875,
46079.19791666666,
46079.20833333334,
46079.21875,
46079.22916666666,
46079.23958333334,
46079.25,
46079.26041666666,
46079.27083333334,
46079.28125,
46079.29166666666,
46079.30208333334,
46079.3125,
46079.32291666666,
46079.33333333334,
46079.34375,
46079.35416666666,
46079.36458333334,
46079.375,
46079.38541666666,
46079.39583333334,
46079.40625,
46079.41666666666,
46079.42708333334,
46079.4375,
46079.44791666666,
46079.45833333334,
46079.46875,
46079.47916666666,
46079.48958333334,
46079.5,
46079.51041666666,
46079.52083333334,
46079.53125,
46079.54166666666,
46079.55208333334,
46079.5625,
46079.57291666666,
46079.58333333334,
46079.59375,
46079.60416666666,
46079.61458333334,
46079.625,
46079.63541666666,
46079.64583333334,
46079.65625,
46079.66666666666,
46079.67708333334,
46079.6875,
46079.69791666666,
46079.70833333334,
46079.71875,
46079.72916666666,
46079.73958333334,
46079.75,
46079.76041666666,
46079.77083333334,
46079.78125,
46079.79166666666,
46079.80208333334,
46079.8125,
46079.82291666666,
46079.83333333334,
46079.84375,
46079.85416666666,
46079.86458333334,
46079.875,
46079.88541666666,
46079.89583333334,
46079.90625,
46079.91666666666,
46079.92708333334,
46079.9375,
46079.94791666666,
46079.95833333334,
46079.96875,
46079.97916666666,
46079.98958333334,
46080
)

# New Notified Production (MW) (column B) values for rows 2..97
$bValues = @(
0,
0,
0,
0,
0,
0,
0,
0,
0,
0,
0,
0,
0,
0,
0,
0,
0,
0,
0,
0,
0.821,
0.881,
0.998,
1.484,
30.437,
61.975,
111.379,
178.261,
408.492,
546.2910000000001,
705.248,
848.6799999999999,
1114.356,
1233.234,
1386.799,
1525.191,
1753.672,
1867.273,
1963.695,
2053.522,
2154.198,
2219.077,
2274.671,
2325.954,
2373.342,
2394.107,
2398.982,
2390.419,
2371.276,
2351.192,
2319.383,
2271.275,
2211.691,
2142.674,
2045.948,
1960.052,
1815.654,
1688.734,
1551.178,
1408.017,
1072.701,
974.074,
793.741,
628.609,
353.864,
236.208,
136.496,
66.14,
12.678,
10.714,
9.949,
9.840999999999999,
6.972,
5.114,
5.079,
5.093,
5.071,
4.731,
0,
2.331,
0.531,
0,
0,
0,
0,
0,
0,
0,
0,
0,
0,
0,
0,
0,
0,
0
)

$n = 96
$aArr = New-Object "object[,]" $n,1
$bArr = New-Object "object[,]" $n,1
for ($i = 0; $i -lt $n; $i++) {
    $aArr[$i,0] = $aValues[$i]
    $bArr[$i,0] = $bValues[$i]
}

$ws.Range("A2:A97").Value = $aArr
$ws.Range("B2:B97").Value = $bArr

Write-Host "Updated rows 2-97 with latest GESS model data."